$d = $word.ActiveDocument

$replacements = @(
    @("2023-11-15 Wednesday", "2023-11-16 Thursday"),
    @("76×28=", "13×98="),
    @("82×62=", "47×24="),
    @("66×39=", "71×60="),
    @("64×82=", "30×72="),
    @("15×18=", "15×31="),
    @("17×38=", "20×24="),
    @("94×38=", "65×68="),
    @("31×28=", "56×48="),
    @("16×77=", "82×79="),
    @("73×19=", "35×92="),
    @("40×28=", "95×80="),
    @("40×68=", "81×34="),
    @("58×60=", "40×79="),
    @("32×52=", "49×31="),
    @("36×14=", "76×14="),
    @("19×67=", "12×14="),
    @("62×80=", "14×27="),
    @("16×94=", "74×38="),
    @("61×68=", "37×63="),
    @("15×76=", "71×52="),
    @("39×98=", "91×40="),
    @("63×29=", "27×97="),
    @("81×92=", "95×82="),
    @("92×23=", "74×26="),
    @("62×55=", "70×34=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
